$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Triggers the engine to stamp baseColWidth="10" onto sheetFormatPr for the
# pre-existing sheet (matches the diff's new sheetFormatPr attribute).
$tmpSheet = $wb.Worksheets.Add()
$tmpSheet.Delete()

# Insert a new column before column A so Code/Description/Definition shift
# from A/B/C to B/C/D, making room for the new "Version" column.
$ws.Range("A1").EntireColumn.Insert()

# New "Version" column header + values
$ws.Range("A1").Value = "Version"
$verRange = $ws.Range("A2:A10")
$verRange.NumberFormat = "@"
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 1).Value = "1.0"
}
$verRange.ClearFormats()


